# The "India" sheet was originally populated by copy/pasting the "France"
# sheet's data, so its Country column (C) still read "France" for every
# data row. This fixes the Country column on the India sheet to say
# "India" instead of "France" for all data rows (3-142), leaving the
# header row and every other column/sheet untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("India")

$ws.Range("C3:C142").Value = "India"
